{"js": "// Helper: replace the unique occurrence of `oldText` with `newText`,\n// preserving the run(s) formatting/markup around it (br, rPr, etc.).\nasync function replaceOnce(oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(oldText) + \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Helper: insert a new italic paragraph right after the paragraph whose\n// full text is `afterText`, with the body `newText`.\nasync function insertItalicParagraphAfter(afterText, newText) {\n  const results = context.document.body.search(afterText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(afterText) + \" but found \" + results.items.length\n    );\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  const newPara = para.insertParagraph(newText, \"After\");\n  newPara.font.set({ italic: true });\n  await context.sync();\n}\n\n// 1) Heading3: \"Advanced treatment and reuse water\" -> \"Water Advanced treatment\"\nawait replaceOnce(\"Advanced treatment and reuse water\", \"Water Advanced treatment\");\n\n// 2) \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2025\"\nawait replaceOnce(\"Ativa\u00e7\u00e3o: 01/01/2012\", \"Ativa\u00e7\u00e3o: 01/01/2025\");\n\n// 3) Insert English translation paragraph after the \"Objetivos\" body paragraph.\nawait insertItalicParagraphAfter(\n  \"Fornecer aos alunos capacidade para: Decidir sobre a melhor alternativa a ser adotada para tratamento avan\u00e7ado, em fun\u00e7\u00e3o da qualidade do efluente a tratar e dos objetivos do reuso ou limita\u00e7\u00f5es do corpo receptor. Conhecer as tecnologias dispon\u00edveis para tratamento avan\u00e7ado de \u00e1guas residu\u00e1rias. Dispor de conceitos e conhecimentos para pr\u00e9-dimensionamento e especifica\u00e7\u00e3o de equipamentos para sistemas de tratamento avan\u00e7ado.\",\n  \"Provide students with the ability to: Decide on the best alternative to be adopted for advanced treatment, depending on the quality of the effluent to be treated and the reuse objectives or limitations of the receiving body. Know the technologies available for advanced wastewater treatment. Have concepts and knowledge for pre-sizing and specification of equipment for advanced treatment systems.\"\n);\n\n// 4) \"Programa resumido\" body: replace text, then add italic translation paragraph.\nawait replaceOnce(\n  \"Tend\u00eancias mundiais no tratamento de \u00e1guas. Processo e opera\u00e7\u00e3o na remo\u00e7\u00e3o de elementos indesej\u00e1veis. Recupera\u00e7\u00e3o de Ambientes Aqu\u00e1ticos.\",\n  \"T\u00e9cnicas alternativas para tratamento. Remo\u00e7\u00e3o de contaminantes importantes. Adequa\u00e7\u00e3o das t\u00e9cnicas avan\u00e7adas \u00e0 exig\u00eancia da legisla\u00e7\u00e3o\"\n);\nawait insertItalicParagraphAfter(\n  \"T\u00e9cnicas alternativas para tratamento. Remo\u00e7\u00e3o de contaminantes importantes. Adequa\u00e7\u00e3o das t\u00e9cnicas avan\u00e7adas \u00e0 exig\u00eancia da legisla\u00e7\u00e3o\",\n  \"Alternative treatment techniques. Removal of important contaminants. Adequacy of advanced techniques to the requirements of legislation\"\n);\n\n// 5) \"Programa\" body: replace text, then add italic translation paragraph.\nawait replaceOnce(\n  \"Tend\u00eancias mundiais sobre tratamento avan\u00e7ado e reuso de \u00e1guas residu\u00e1rias. Determina\u00e7\u00e3o da efici\u00eancia de processos e opera\u00e7\u00f5es em fun\u00e7\u00e3o dos objetivos de reuso da qualidade do afluente a tratar e da obedi\u00eancia a padr\u00f5es de emiss\u00e3o e de qualidade. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o de nitrog\u00eanio e f\u00f3sforo: nitrifica\u00e7\u00e3o, desnitrifica\u00e7\u00e3o, e remo\u00e7\u00e3o qu\u00edmica e biol\u00f3gica de f\u00f3sforo. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, \\\"stripping\\\", coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos. Disposi\u00e7\u00e3o e tratamento de esgotos no solo. Recupera\u00e7\u00e3o de ambientes aqu\u00e1ticos com base na piscicultura e aproveitamento de algas e macr\u00f3fias.\",\n  \"T\u00e9cnicas alterativas para tratamento de \u00e1gua e efluentes. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o dos principais contaminantes indicados na legisla\u00e7\u00e3o vigente. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos.\"\n);\nawait insertItalicParagraphAfter(\n  \"T\u00e9cnicas alterativas para tratamento de \u00e1gua e efluentes. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o dos principais contaminantes indicados na legisla\u00e7\u00e3o vigente. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos.\",\n  \"Alternative techniques for water and effluent treatment. Processes and operations applied to the removal of the main contaminants indicated in current legislation. Operations and processes for removing specific contaminants: adsorption on activated carbon, chemical oxidation, coagulation-flocculation (sedimentation and flotation), ion exchange, reverse osmosis, membrane filtration, filtration in porous media.\"\n);\n\n// 6) \"M\u00e9todo:\" run text.\nawait replaceOnce(\n  \"Ser\u00e3o ministradas aulas expositivas convencionais, associadas \u00e0 exposi\u00e7\u00e3o de v\u00eddeos e slides sobre sistemas de tratamento avan\u00e7ado. Al\u00e9m disso ser\u00e3o efetuadas visitas e ser\u00e3o desenvolvidos exerc\u00edcios orientados.\",\n  \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n);\n\n// 7) \"Crit\u00e9rio:\" run text.\nawait replaceOnce(\n  \"Ser\u00e3o aplicadas duas provas (1o. e 2o. bimestres), com peso 8,0 e tamb\u00e9m ser\u00e1 entregue lista de exerc\u00edcios, com peso 2,0.\",\n  \"M\u00e9dia ponderada das notas atribu\u00eddas \u00e0s provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n);\n\n// 8) \"Norma de recupera\u00e7\u00e3o:\" run text.\nawait replaceOnce(\n  \"Elabora\u00e7\u00e3o de monografia, com tema escolhido pelo docente, enfocando mat\u00e9ria em que o aluno demonstrou menor habilidade (peso: 3,0); e prova escrita sobre todfa a mat\u00e9ria da disciplina (peso: 7,0).\",\n  \"1 (uma) prova de recupera\u00e7\u00e3o (R), sendo considerado aprovado se R >= 5,0.\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Once($doc, $oldText, $newText) {\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found for replace: $oldText\"\n    }\n}\n\nfunction Find-ParagraphIndex($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        # Strip trailing paragraph/cell mark character(s) before comparing.\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Insert-ItalicParagraphAfter($doc, $afterText, $newText) {\n    $idx = Find-ParagraphIndex $doc $afterText\n    if ($idx -eq -1) {\n        throw \"Paragraph not found: $afterText\"\n    }\n    $p = $doc.Paragraphs.Item($idx)\n    [void]$p.Range.InsertParagraphAfter()\n    $newPara = $doc.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = $newText\n    # Re-fetch and trim the trailing paragraph mark off the range before\n    # italicising, otherwise Word also stamps the paragraph-mark rPr.\n    $rng = $doc.Paragraphs.Item($idx + 1).Range\n    $rng.MoveEnd(1, -1)\n    $rng.Font.Italic = $true\n}\n\n# 1) Heading3: \"Advanced treatment and reuse water\" -> \"Water Advanced treatment\"\nReplace-Once $d \"Advanced treatment and reuse water\" \"Water Advanced treatment\"\n\n# 2) \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2025\"\nReplace-Once $d \"Ativa\u00e7\u00e3o: 01/01/2012\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\n# 3) Insert English translation paragraph after the \"Objetivos\" body paragraph.\nInsert-ItalicParagraphAfter $d `\n    \"Fornecer aos alunos capacidade para: Decidir sobre a melhor alternativa a ser adotada para tratamento avan\u00e7ado, em fun\u00e7\u00e3o da qualidade do efluente a tratar e dos objetivos do reuso ou limita\u00e7\u00f5es do corpo receptor. Conhecer as tecnologias dispon\u00edveis para tratamento avan\u00e7ado de \u00e1guas residu\u00e1rias. Dispor de conceitos e conhecimentos para pr\u00e9-dimensionamento e especifica\u00e7\u00e3o de equipamentos para sistemas de tratamento avan\u00e7ado.\" `\n    \"Provide students with the ability to: Decide on the best alternative to be adopted for advanced treatment, depending on the quality of the effluent to be treated and the reuse objectives or limitations of the receiving body. Know the technologies available for advanced wastewater treatment. Have concepts and knowledge for pre-sizing and specification of equipment for advanced treatment systems.\"\n\n# 4) \"Programa resumido\" body: replace text, then add italic translation paragraph.\nReplace-Once $d `\n    \"Tend\u00eancias mundiais no tratamento de \u00e1guas. Processo e opera\u00e7\u00e3o na remo\u00e7\u00e3o de elementos indesej\u00e1veis. Recupera\u00e7\u00e3o de Ambientes Aqu\u00e1ticos.\" `\n    \"T\u00e9cnicas alternativas para tratamento. Remo\u00e7\u00e3o de contaminantes importantes. Adequa\u00e7\u00e3o das t\u00e9cnicas avan\u00e7adas \u00e0 exig\u00eancia da legisla\u00e7\u00e3o\"\n\nInsert-ItalicParagraphAfter $d `\n    \"T\u00e9cnicas alternativas para tratamento. Remo\u00e7\u00e3o de contaminantes importantes. Adequa\u00e7\u00e3o das t\u00e9cnicas avan\u00e7adas \u00e0 exig\u00eancia da legisla\u00e7\u00e3o\" `\n    \"Alternative treatment techniques. Removal of important contaminants. Adequacy of advanced techniques to the requirements of legislation\"\n\n# 5) \"Programa\" body: replace text, then add italic translation paragraph.\nReplace-Once $d `\n    'Tend\u00eancias mundiais sobre tratamento avan\u00e7ado e reuso de \u00e1guas residu\u00e1rias. Determina\u00e7\u00e3o da efici\u00eancia de processos e opera\u00e7\u00f5es em fun\u00e7\u00e3o dos objetivos de reuso da qualidade do afluente a tratar e da obedi\u00eancia a padr\u00f5es de emiss\u00e3o e de qualidade. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o de nitrog\u00eanio e f\u00f3sforo: nitrifica\u00e7\u00e3o, desnitrifica\u00e7\u00e3o, e remo\u00e7\u00e3o qu\u00edmica e biol\u00f3gica de f\u00f3sforo. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, \"stripping\", coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos. Disposi\u00e7\u00e3o e tratamento de esgotos no solo. Recupera\u00e7\u00e3o de ambientes aqu\u00e1ticos com base na piscicultura e aproveitamento de algas e macr\u00f3fias.' `\n    \"T\u00e9cnicas alterativas para tratamento de \u00e1gua e efluentes. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o dos principais contaminantes indicados na legisla\u00e7\u00e3o vigente. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos.\"\n\nInsert-ItalicParagraphAfter $d `\n    \"T\u00e9cnicas alterativas para tratamento de \u00e1gua e efluentes. Processos e opera\u00e7\u00f5es aplicadas \u00e0 remo\u00e7\u00e3o dos principais contaminantes indicados na legisla\u00e7\u00e3o vigente. Opera\u00e7\u00f5es e processos para remo\u00e7\u00e3o de contaminantes espec\u00edficos: adsor\u00e7\u00e3o em carv\u00e3o ativado, oxida\u00e7\u00e3o qu\u00edmica, coagula\u00e7\u00e3o-flocula\u00e7\u00e3o (sedimenta\u00e7\u00e3o e flota\u00e7\u00e3o), troca i\u00f4nica, osmose reversa, filtra\u00e7\u00e3o em membranas, filtra\u00e7\u00e3o em meios porosos.\" `\n    \"Alternative techniques for water and effluent treatment. Processes and operations applied to the removal of the main contaminants indicated in current legislation. Operations and processes for removing specific contaminants: adsorption on activated carbon, chemical oxidation, coagulation-flocculation (sedimentation and flotation), ion exchange, reverse osmosis, membrane filtration, filtration in porous media.\"\n\n# 6) \"M\u00e9todo:\" run text.\nReplace-Once $d `\n    \"Ser\u00e3o ministradas aulas expositivas convencionais, associadas \u00e0 exposi\u00e7\u00e3o de v\u00eddeos e slides sobre sistemas de tratamento avan\u00e7ado. Al\u00e9m disso ser\u00e3o efetuadas visitas e ser\u00e3o desenvolvidos exerc\u00edcios orientados.\" `\n    \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n\n# 7) \"Crit\u00e9rio:\" run text.\nReplace-Once $d `\n    \"Ser\u00e3o aplicadas duas provas (1o. e 2o. bimestres), com peso 8,0 e tamb\u00e9m ser\u00e1 entregue lista de exerc\u00edcios, com peso 2,0.\" `\n    \"M\u00e9dia ponderada das notas atribu\u00eddas \u00e0s provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n\n# 8) \"Norma de recupera\u00e7\u00e3o:\" run text.\nReplace-Once $d `\n    \"Elabora\u00e7\u00e3o de monografia, com tema escolhido pelo docente, enfocando mat\u00e9ria em que o aluno demonstrou menor habilidade (peso: 3,0); e prova escrita sobre todfa a mat\u00e9ria da disciplina (peso: 7,0).\" `\n    \"1 (uma) prova de recupera\u00e7\u00e3o (R), sendo considerado aprovado se R >= 5,0.\"\n"}
